$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had two "header" rows ("situação do domicílio" at row 5 and
# "grandes regiões" at row 8) that carried no data. The correction removes
# those two label rows entirely (and their now-unused shared strings) and
# shifts the remaining data rows up so the data that used to sit one/two
# rows below each label now lines up with the row above it.

# Remove the "situação do domicílio" row (row 5) - everything below shifts up.
$ws.Rows(5).Delete()

# After the first deletion, "grandes regiões" (originally row 8) is now row 7.
$ws.Rows(7).Delete()
